$wb = $excel.ActiveWorkbook

$oldId = "2f129653-0815-4eb5-95e8-1d677fdf504c"
$newId = "0301ae3f-9670-4e0c-9f46-a0e5c0433730"

$oldHash = "352a1d77f4003e0b4581c51b35e5ee52e49012bd"
$newHash = "f7de731860baa984aab37f52c53b2e5acc9ba6e0"

# The hyperlink target (rels) URL itself is left untouched by the original
# edit -- it still points at the OLD id; only the display text changes.
$linkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3ccb43bd0782758c4d1976a1fbc6a737b58bb185/e2e/$oldId.md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---- Overview sheet ----
# A2: File Name
$wsOverview.Range("A2").Value = "$newId.md"
# B2: Path And Name (hyperlink cell) - update both the cell text and the
# hyperlink's display text, keep its target address unchanged.
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkTarget, $null, $null, "e2e\$newId.md") | Out-Null
# G2: Latest HO Xliff Generate Date
$wsOverview.Range("G2").Value = "2016-11-14 06:44:56"

# ---- zh-cn sheet ----
# A2: Source File Name (hyperlink cell)
$wsZh.Range("A2").Value = "$newId.md"
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $linkTarget, $null, $null, "$newId.md") | Out-Null
# G2: Latest Handoff File
$wsZh.Range("G2").Value = "$newId.$newHash.zh-cn.xlf"
# H2: Latest Handoff Datetime
$wsZh.Range("H2").Value = "2016-11-14 06:44:43"

# ---- de-de sheet ----
# A2: Source File Name (hyperlink cell)
$wsDe.Range("A2").Value = "$newId.md"
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $linkTarget, $null, $null, "$newId.md") | Out-Null
# G2: Latest Handoff File
$wsDe.Range("G2").Value = "$newId.$newHash.de-de.xlf"
